# Add a new "Base-Model" indicator column right after the Variable Name
# column (i.e. insert a new column B), flag the variables that belong to
# the base model with "Y", and leave the rest of that column blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, empty column before the current column B. This shifts the
# old Description/Range/Type columns (and every data column after them)
# one slot to the right, exactly like Excel's own "Insert Column".
$ws.Columns("B:B").Insert()

# Header for the freshly inserted column.
$ws.Cells.Item(1, 2).Value = "Base-Model"

# Rows (by their 1-based worksheet row number) whose variable is part of
# the base model -- mark them with "Y" in the new column.
$baseModelRows = @(12, 13, 14, 15, 32, 33, 37, 48, 54, 55)
foreach ($r in $baseModelRows) {
    $ws.Cells.Item($r, 2).Value = "Y"
}

# Size the new column to fit its (short) contents, closely matching the
# width Excel itself would have picked when autofitting this column.
$ws.Columns("B:B").ColumnWidth = 9.8

# Match the author's final cursor position/selection in the sheet.
$ws.Range("E13").Select()
